$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Coin/Link/Price/Volume columns so Excel does not
# auto-convert numeric-looking strings like "1.002" or dotted-thousands
# strings like "23.848.86" into actual numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.848.86'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '1.640.90'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '309.83'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").Value = '0.3878'
$ws.Range("E7").Value = '  -0.87%  '
$ws.Range("D8").Value = '0.3829'
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").Value = '50.48'
$ws.Range("E9").Value = '  -2.38%  '
$ws.Range("D10").Value = '1.325'
$ws.Range("E10").Value = '  -3.74%  '
$ws.Range("D11").Value = '1.002'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '0.08389'
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("D13").Value = '23.68'
$ws.Range("E13").Value = '  -2.36%  '
$ws.Range("D14").Value = '6.972'
$ws.Range("E14").Value = '  -4.13%  '
$ws.Range("D15").Value = '7.846'
$ws.Range("E15").Value = '  -3.83%  '
$ws.Range("D16").Value = '0.00001309'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '1.642.46'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '93.65'
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("D19").Value = '0.06956'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '19.41'
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("D21").Value = '6.872'
$ws.Range("E21").Value = '  -1.57%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = '13.58'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("D24").Value = '23.856.52'
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").Value = '2.442'
$ws.Range("E25").Value = '  -2.52%  '
$ws.Range("D26").Value = '2.882'
$ws.Range("E26").Value = '  -9.00%  '
$ws.Range("D27").Value = '21.81'
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '153.08'
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("B29").Value = 'HuobiToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D29").Value = '5.572'
$ws.Range("E29").Value = '  +4.98%  '
$ws.Range("D30").Value = '136.69'
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '7.641'
$ws.Range("E31").Value = '  -3.19%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.495'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = '1.820.13'
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("D34").Value = '0.07994'
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").Value = '0.9788'
$ws.Range("E35").Value = '  -6.77%  '
$ws.Range("D36").Value = '0.02894'
$ws.Range("E36").Value = '  -4.78%  '
$ws.Range("D37").Value = '6.571'
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").Value = '0.2653'
$ws.Range("E38").Value = '  -2.47%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '10.38'
$ws.Range("E39").Value = '  -7.54%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.09080'
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").Value = '0.7485'
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").Value = '13.28'
$ws.Range("E42").Value = '  -1.93%  '
$ws.Range("D43").Value = '1.413'
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").Value = '16.48'
$ws.Range("E44").Value = '  -0.53%  '
$ws.Range("D45").Value = '0.6883'
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("D46").Value = '2.416'
$ws.Range("E46").Value = '  -3.76%  '
$ws.Range("D47").Value = '4.084'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("D49").Value = '0.08217'
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").Value = '134.01'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = '1.213'
$ws.Range("E51").Value = '  -2.57%  '

Write-Output "Updated $([int]109) cells"
